# Scheduled-runner price/profit refresh across the Leve-crafting tables
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each block below rewrites the
# market-price-derived columns (H..N) for the specific rows whose source
# prices changed; unaffected columns/rows are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 73166.36
$ws.Range("I28").Value = 144376.86
$ws.Range("K28").Value = 144376.86
$ws.Range("M28").Value = -143891.86

$ws.Range("H33").Value = 244.44444
$ws.Range("I33").Value = 244.44444
$ws.Range("K33").Value = 244.44444
$ws.Range("M33").Value = -15.44443999999999

$ws.Range("H82").Value = 1679.5
$ws.Range("I82").Value = 1679.5
$ws.Range("K82").Value = 5038.5
$ws.Range("M82").Value = -4632.5

$ws.Range("H85").Value = 1679.5
$ws.Range("I85").Value = 1679.5
$ws.Range("K85").Value = 5038.5
$ws.Range("M85").Value = -3634.5

$ws.Range("H98").Value = 2679.52
$ws.Range("J98").Value = 4569.091
$ws.Range("L98").Value = 4569.091
$ws.Range("N98").Value = -7565.091

$ws.Range("H122").Value = 2679.52
$ws.Range("J122").Value = 4569.091
$ws.Range("L122").Value = 13707.273
$ws.Range("N122").Value = -18607.273

$ws.Range("H137").Value = 3395.8147
$ws.Range("I137").Value = 2121.6858
$ws.Range("K137").Value = 6365.057400000001
$ws.Range("M137").Value = -3815.057400000001

$ws.Range("H138").Value = 4275.787
$ws.Range("I138").Value = 2681.4614
$ws.Range("K138").Value = 8044.3842
$ws.Range("M138").Value = -2904.3842

$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws.Range("H141").Value = 2839.7693
$ws.Range("I141").Value = 1441.7
$ws.Range("K141").Value = 4325.1
$ws.Range("M141").Value = 854.8999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5977.5864
$ws.Range("I2").Value = 5748.222
$ws.Range("J2").Value = 6352.909
$ws.Range("K2").Value = 5748.222
$ws.Range("L2").Value = 6352.909
$ws.Range("M2").Value = -5635.222
$ws.Range("N2").Value = -6578.909

$ws.Range("H4").Value = 142857420
$ws.Range("I4").Value = 299.66666
$ws.Range("J4").Value = 250000260
$ws.Range("K4").Value = 299.66666
$ws.Range("L4").Value = 250000260
$ws.Range("M4").Value = -183.66666
$ws.Range("N4").Value = -250000492

$ws.Range("H45").Value = 4754.95
$ws.Range("J45").Value = 5582.25
$ws.Range("L45").Value = 5582.25
$ws.Range("N45").Value = -6336.25

$ws.Range("H97").Value = 2282.375
$ws.Range("I97").Value = 2393.7144
$ws.Range("J97").Value = 1503
$ws.Range("K97").Value = 2393.7144
$ws.Range("L97").Value = 1503
$ws.Range("M97").Value = -1897.7144
$ws.Range("N97").Value = -2495

$ws.Range("H102").Value = 1989.24
$ws.Range("I102").Value = 1910.0869
$ws.Range("K102").Value = 1910.0869
$ws.Range("M102").Value = -288.0869

$ws.Range("H110").Value = 1206
$ws.Range("I110").Value = 1305.2222
$ws.Range("J110").Value = 1027.4
$ws.Range("K110").Value = 1305.2222
$ws.Range("L110").Value = 1027.4
$ws.Range("M110").Value = 739.7778000000001
$ws.Range("N110").Value = -5117.4

$ws.Range("H116").Value = 5977.5864
$ws.Range("I116").Value = 5748.222
$ws.Range("J116").Value = 6352.909
$ws.Range("K116").Value = 5748.222
$ws.Range("L116").Value = 6352.909
$ws.Range("M116").Value = -3454.222
$ws.Range("N116").Value = -10940.909

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5977.5864
$ws.Range("I3").Value = 5748.222
$ws.Range("J3").Value = 6352.909
$ws.Range("K3").Value = 5748.222
$ws.Range("L3").Value = 6352.909
$ws.Range("M3").Value = -5634.222
$ws.Range("N3").Value = -6580.909

$ws.Range("H132").Value = 49999.8
$ws.Range("J132").Value = 49999.8
$ws.Range("L132").Value = 49999.8
$ws.Range("N132").Value = -60119.8

$ws.Range("H139").Value = 99740
$ws.Range("J139").Value = 99740
$ws.Range("L139").Value = 99740
$ws.Range("N139").Value = -110020

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 445.5
$ws.Range("J22").Value = 550
$ws.Range("L22").Value = 550
$ws.Range("N22").Value = -1250

$ws.Range("H58").Value = 481338.62
$ws.Range("I58").Value = 1002230.1
$ws.Range("J58").Value = 7800.909
$ws.Range("K58").Value = 1002230.1
$ws.Range("L58").Value = 7800.909
$ws.Range("M58").Value = -1002027.1
$ws.Range("N58").Value = -8206.909

$ws.Range("H94").Value = 2911
$ws.Range("J94").Value = 3262.8333
$ws.Range("L94").Value = 3262.8333
$ws.Range("N94").Value = -4164.8333

$ws.Range("H132").Value = 5059.048
$ws.Range("I132").Value = 4527.2256
$ws.Range("K132").Value = 13581.6768
$ws.Range("M132").Value = -11051.6768

$ws.Range("H136").Value = 481338.62
$ws.Range("I136").Value = 1002230.1
$ws.Range("J136").Value = 7800.909
$ws.Range("K136").Value = 3006690.3
$ws.Range("L136").Value = 23402.727
$ws.Range("M136").Value = -3004140.3
$ws.Range("N136").Value = -28502.727

$ws.Range("H139").Value = 94053.336
$ws.Range("J139").Value = 93580
$ws.Range("L139").Value = 93580
$ws.Range("N139").Value = -103860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H51").Value = 2625
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = -8540
$ws.Range("N51").Value = -8420

$ws.Range("H131").Value = 5143.857
$ws.Range("I131").Value = 1116.375
$ws.Range("J131").Value = 7622.3076
$ws.Range("K131").Value = 3349.125
$ws.Range("L131").Value = 22866.9228
$ws.Range("M131").Value = 1690.875
$ws.Range("N131").Value = -32946.9228

$ws.Range("H134").Value = 1643.5
$ws.Range("I134").Value = 1643.5
$ws.Range("K134").Value = 4930.5
$ws.Range("M134").Value = 139.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8048
$ws.Range("I122").Value = 7576.316
$ws.Range("J122").Value = 9328.286
$ws.Range("K122").Value = 22728.948
$ws.Range("L122").Value = 27984.858
$ws.Range("M122").Value = -20278.948
$ws.Range("N122").Value = -32884.858

$ws.Range("H126").Value = 3286.9333
$ws.Range("I126").Value = 3102.5
$ws.Range("K126").Value = 9307.5
$ws.Range("M126").Value = -6837.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11797.333
$ws.Range("I41").Value = 11578
$ws.Range("J41").Value = 11841.2
$ws.Range("K41").Value = 11578
$ws.Range("L41").Value = 11841.2
$ws.Range("M41").Value = -11188
$ws.Range("N41").Value = -12621.2

$ws.Range("H81").Value = 11019.23
$ws.Range("I81").Value = 1128.6666
$ws.Range("J81").Value = 19496.857
$ws.Range("K81").Value = 2257.3332
$ws.Range("L81").Value = 38993.714
$ws.Range("M81").Value = -1196.3332
$ws.Range("N81").Value = -41115.714

$ws.Range("H84").Value = 11019.23
$ws.Range("I84").Value = 1128.6666
$ws.Range("J84").Value = 19496.857
$ws.Range("K84").Value = 11286.666
$ws.Range("L84").Value = 194968.57
$ws.Range("M84").Value = -5982.666000000001
$ws.Range("N84").Value = -205576.57

$ws.Range("H97").Value = 21853.334
$ws.Range("J97").Value = 21853.334
$ws.Range("L97").Value = 21853.334
$ws.Range("N97").Value = -23835.334
